# Regione: REGIONE_CAMPANIA, ASL: ASL_NA_3_SUD, Issuer: integrity:S1#VICAMPANIA3SUD
#
# Row 2 (the single data row) gets a new "REGIONE_LAZIO" value in column B
# (previously empty) and the IHE/workflow identifiers in D, E and F are
# refreshed to a new gateway run (new workflowInstanceId suffix, new
# UAT_GTW_ID and a new event timestamp). Column C (the patient id) is
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the existing identifiers/timestamp in place first ...
$ws.Range("D2").Value = "2.16.840.1.113883.2.9.2.120.4.4.b0f3ffcf25ce2aafc7dc901e2febc51f43837f4ca0fe3b6d1b02194e9047b6db.f6fb3ca4ab^^^^urn:ihe:iti:xdw:2013:workflowInstanceId"
$ws.Range("E2").Value = "2.16.840.1.113883.2.9.2.110.4.4^UAT_GTW_ID1721654247090"
$ws.Range("F2").Value = "22-07-2024:15:17:28"

# ... then populate the previously-empty B2 cell with the new region label.
$ws.Range("B2").Value = "REGIONE_LAZIO"
